$wb = $excel.ActiveWorkbook

# zh-cn sheet: update "Correspond Handoff Datetime" (D4) and
# "Correspond Handback DateTime" (G4) for the 1e64d118... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-17 10:07:32"
$wsZhCn.Range("G4").Value = "2016-01-17 10:08:18"

# de-de sheet: update "Correspond Handoff Datetime" (D4) and
# "Correspond Handback DateTime" (G4) for the 1e64d118... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-17 10:07:43"
$wsDeDe.Range("G4").Value = "2016-01-17 10:08:37"
